$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "F2", "F2rl2", "ECs", 2, 1, 0.627495, 1.25499, 0.4221230570832646, 0.3536267886120559, 1, 0.5, 0.2245445, 0.449089, 0.03525495836513261, 0.0237827926988904, 0.1409005510275, 0.5636022041100001, 0.01488193080243299, 0.008410232606334861),
    @("ECs", "F2", "F2rl2", "Inflammatory-Mac", 2, 1, 0.627495, 1.25499, 0.4221230570832646, 0.3536267886120559, 3, 1, 1.393304333333333, 4.179913, 0.2187579132939151, 0.2213592503454706, 0.874291502645, 5.24574901587, 0.09234275912078314, 0.0782785608292409),
    @("ECs", "F2", "F2rl2", "Neutrophils", 2, 1, 0.627495, 1.25499, 0.4221230570832646, 0.3536267886120559, 3, 1, 2.796960666666667, 8.390882000000001, 0.4391411584440807, 0.4443631600124939, 1.75507883353, 10.53047300118, 0.1853716082935016, 0.1571387172527234),
    @("ECs", "F2", "F2rl2", "Resolving-Mac", 2, 1, 0.627495, 1.25499, 0.4221230570832646, 0.3536267886120559, 3, 1, 1.954351333333333, 5.863053999999999, 0.3068459698968715, 0.3104947969431451, 1.22634568991, 7.358074139459999, 0.1295267588665468, 0.1097992779237568),
    @("FAPs", "F2", "F2rl2", "ECs", 3, 1, 0.524746, 1.574238, 0.3530026306380366, 0.4435833978366884, 1, 0.5, 0.2245445, 0.449089, 0.03525495836513261, 0.0237827926988904, 0.117828828197, 0.706972969182, 0.01244509304592627, 0.01054965199541939),
    @("FAPs", "F2", "F2rl2", "Inflammatory-Mac", 3, 1, 0.524746, 1.574238, 0.3530026306380366, 0.4435833978366884, 3, 1, 1.393304333333333, 4.179913, 0.2187579132939151, 0.2213592503454706, 0.7311308756993333, 6.580177881294, 0.07722211886563954, 0.098191288410826),
    @("FAPs", "F2", "F2rl2", "Neutrophils", 3, 1, 0.524746, 1.574238, 0.3530026306380366, 0.4435833978366884, 3, 1, 2.796960666666667, 8.390882000000001, 0.4391411584440807, 0.4443631600124939, 1.467693921990667, 13.209245297916, 0.1550179841521953, 0.1971121203917901),
    @("FAPs", "F2", "F2rl2", "Resolving-Mac", 3, 1, 0.524746, 1.574238, 0.3530026306380366, 0.4435833978366884, 3, 1, 1.954351333333333, 5.863053999999999, 0.3068459698968715, 0.3104947969431451, 1.025538044761333, 9.229842402851999, 0.1083174345742754, 0.137730337038653),
    @("MuSCs", "F2", "F2rl2", "ECs", 2, 1, 0.2831585, 0.566317, 0.1904839586914821, 0.1595748667689891, 1, 0.5, 0.2245445, 0.449089, 0.03525495836513261, 0.0237827926988904, 0.06358168380324999, 0.254326735213, 0.006715504032893841, 0.003795135976319922),
    @("MuSCs", "F2", "F2rl2", "Inflammatory-Mac", 2, 1, 0.2831585, 0.566317, 0.1904839586914821, 0.1595748667689891, 3, 1, 1.393304333333333, 4.179913, 0.2187579132939151, 0.2213592503454706, 0.3945259650701666, 2.367155790421, 0.04166987331931294, 0.03532337288196178),
    @("MuSCs", "F2", "F2rl2", "Neutrophils", 2, 1, 0.2831585, 0.566317, 0.1904839586914821, 0.1595748667689891, 3, 1, 2.796960666666667, 8.390882000000001, 0.4391411584440807, 0.4443631600124939, 0.7919831869323334, 4.751899121594001, 0.08364934628479187, 0.07090919205604071),
    @("MuSCs", "F2", "F2rl2", "Resolving-Mac", 2, 1, 0.2831585, 0.566317, 0.1904839586914821, 0.1595748667689891, 3, 1, 1.954351333333333, 5.863053999999999, 0.3068459698968715, 0.3104947969431451, 0.5533911920196666, 3.320347152117999, 0.05844923505448343, 0.04954716585466671),
    @("Resolving-Mac", "F2", "F2rl2", "ECs", 1, 0.3333333333333333, 0.051122, 0.153366, 0.03439035358721687, 0.04321494678226644, 1, 0.5, 0.2245445, 0.449089, 0.03525495836513261, 0.0237827926988904, 0.011479163929, 0.068874983574, 0.00121243048387952, 0.001027772120816223),
    @("Resolving-Mac", "F2", "F2rl2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.051122, 0.153366, 0.03439035358721687, 0.04321494678226644, 3, 1, 1.393304333333333, 4.179913, 0.2187579132939151, 0.2213592503454706, 0.07122850412866666, 0.641056537158, 0.007523161988179469, 0.009566028223441906),
    @("Resolving-Mac", "F2", "F2rl2", "Neutrophils", 1, 0.3333333333333333, 0.051122, 0.153366, 0.03439035358721687, 0.04321494678226644, 3, 1, 2.796960666666667, 8.390882000000001, 0.4391411584440807, 0.4443631600124939, 0.1429862232013334, 1.286876008812, 0.01510221971359197, 0.01920313031193967),
    @("Resolving-Mac", "F2", "F2rl2", "Resolving-Mac", 1, 0.3333333333333333, 0.051122, 0.153366, 0.03439035358721687, 0.04321494678226644, 3, 1, 1.954351333333333, 5.863053999999999, 0.3068459698968715, 0.3104947969431451, 0.09991034886266666, 0.8991931397639998, 0.01055254140156591, 0.01341801612606864)
)

$numRows = $data.Count
$numCols = $data[0].Count

# Delete the rows that are being removed (old rows 18 and 19), since the
# new data only spans rows 2..17 (16 data rows instead of 18).
$ws.Rows.Item(18).Resize(2).EntireRow.Delete() | Out-Null

for ($r = 0; $r -lt $numRows; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $numCols; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowValues[$c]
    }
}
